$d = $word.ActiveDocument

# Replace the entire body content with a single paragraph of large,
# bold, red, yellow-highlighted text, wrapped in grammar-check markers,
# while leaving the section properties (page size/margins/etc.) intact.
$rpr = '<w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="96"/><w:szCs w:val="96"/>'
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:rPr>' + $rpr + '</w:rPr></w:pPr>' +
       '<w:proofErr w:type="gramStart"/>' +
       '<w:r><w:rPr>' + $rpr + '<w:highlight w:val="yellow"/></w:rPr><w:t>MAKING  changes</w:t></w:r>' +
       '<w:proofErr w:type="gramEnd"/>' +
       '</w:p>' +
       '<w:sectPr>' +
       '<w:pgSz w:w="12240" w:h="15840"/>' +
       '<w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/>' +
       '<w:cols w:space="720"/>' +
       '<w:docGrid w:linePitch="360"/>' +
       '</w:sectPr>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
